$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15 (shifts existing rows 15-49 down to 16-50)
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new localization key/value
$ws.Range("A15").Value = "BUTTON_CANCEL"
$ws.Range("B15").Value = "Cancel"
$ws.Range("C15").Value = "XXXX"
$ws.Range("D15").Value = "XXXX"
$ws.Range("E15").Value = "XXXX"

# Match the final selection state recorded in the sheet view
$ws.Range("E15").Select()
